$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date (serial 42025 = 2015-01-21) for the two new rows
$newDate = 42025

# Row 7: low priority, DAQLab module
$ws.Range("A7").Value = $newDate
$ws.Range("B7").Value = "low"
$ws.Range("C7").Value = "DAQLab"
$ws.Range("D7").Value = "Develop a uniform naming system for Vchans and modify DLGetUniqueVChanName. The naming convention should be `"<module name>: <task controller name>: <VChanName> <number>`" with the number being optional in case there would be multiple Vchans with the same name. One issue to address in this case is how to change the VChan name if the module or task controller names change"
$ws.Range("E7").Value = "Adrian"
$ws.Rows.Item(7).RowHeight = 60

# Row 8: high priority, DAQLab, task controller module
$ws.Range("A8").Value = $newDate
$ws.Range("B8").Value = "high"
$ws.Range("C8").Value = "DAQLab, task controller"
$ws.Range("D8").Value = "If an error is encountered when executing an UITC and after that the UITC is reset and started again, then the UITC receives a Start event while it is stuck in its Configured state. This is seen when the UITC has child TCs, and it is possible that there is some sort of race condition again between threads. Sometimes an error does not need to occur, it seems that the reset itself is just unreliable."
$ws.Range("E8").Value = "Adrian"
$ws.Rows.Item(8).RowHeight = 60

# Update selection to D8
$ws.Range("D8").Select()
